# Insert a new data row at row 70 (shifting existing rows 70-129 down to 71-130)
# and populate it with the new Albahaca price record for Femacal de La Calera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(70).Insert()

$ws.Cells.Item(70, 1).Value = 3
$ws.Cells.Item(70, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(70, 3).Value = "Coquimbo"
$ws.Cells.Item(70, 4).Value = 44589
$ws.Cells.Item(70, 5).Value = 5
$ws.Cells.Item(70, 6).Value = 100112052
$ws.Cells.Item(70, 7).Value = "Albahaca"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 160
$ws.Cells.Item(70, 11).Value = 4000
$ws.Cells.Item(70, 12).Value = 4500
$ws.Cells.Item(70, 13).Value = 4250
$ws.Cells.Item(70, 14).Value = "$/docena de matas"
$ws.Cells.Item(70, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(70, 16).Value = 708
$ws.Cells.Item(70, 17).Value = 6
$ws.Cells.Item(70, 18).Value = "Hortaliza"
